$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must have their
# NumberFormat forced to text ("@") first, otherwise Excel's COM layer
# auto-converts the assigned string into a floating point number and
# mangles values like '6.33' or '0.259' (e.g. into 6.3300000000000001).
# The source workbook stores every one of these cells as inline text,
# so forcing text formatting keeps the on-disk representation faithful.
$textCoercedCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D31", "D32", "D35", "D39", "D42", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($addr in $textCoercedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (mirrors the authoritative XML diff).
$ws.Range("D2").Value = "25.787.48"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.635.36"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "215.25"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "0.502"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.259"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "0.0638"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("D11").Value = "0.0795"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Value = "4.28"
$ws.Range("E12").Value = "  +0.22%  "
$ws.Range("D13").Value = "1.859.73"
$ws.Range("D14").Value = "1.635.09"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("E15").Value = "  -0.56%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "63.11"
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "25.810.57"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "4.46"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "193.33"
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").Value = "9.97"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "6.33"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("B24").Value = "BinanceUSD"
$ws.Range("C24").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "1.82"
$ws.Range("E25").Value = "  +2.53%  "
$ws.Range("D26").Value = "143.16"
$ws.Range("E26").Value = "  +3.56%  "
$ws.Range("D27").Value = "0.123"
$ws.Range("E27").Value = "  -1.00%  "
$ws.Range("D28").Value = "6.95"
$ws.Range("E28").Value = "  +1.22%  "
$ws.Range("D29").Value = "15.57"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("D31").Value = "0.0493"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("D32").Value = "3.35"
$ws.Range("E32").Value = "  +1.32%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  -0.37%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").Value = "1.135.38"
$ws.Range("E37").Value = "  -0.02%  "
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "0.543"
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  -1.08%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").Value = "5.57"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "100.71"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "0.806"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "1.770.44"
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "55.25"
$ws.Range("E47").Value = "  -0.87%  "
$ws.Range("E48").Value = "  -1.95%  "
$ws.Range("D49").Value = "0.0503"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("D50").Value = "7.49"
$ws.Range("E50").Value = "  -3.30%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.10%  "
